# Updates team matchup-probability matrix values on Sheet1 (North Texas_A)
# with refreshed figures from games pulled March 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.1935483870967742
$ws.Cells.Item(2, 3).Value = 0.5622119815668203
$ws.Cells.Item(2, 10).Value = 0.01382488479262673
$ws.Cells.Item(2, 15).Value = 0.004608294930875576
$ws.Cells.Item(2, 16).Value = 0.1612903225806452
$ws.Cells.Item(2, 19).Value = 0.06451612903225806
$ws.Cells.Item(3, 2).Value = 0.01538461538461539
$ws.Cells.Item(3, 3).Value = 0.06153846153846154
$ws.Cells.Item(3, 10).Value = 0.02307692307692308
$ws.Cells.Item(3, 16).Value = 0.7230769230769231
$ws.Cells.Item(3, 19).Value = 0.1769230769230769
$ws.Cells.Item(4, 10).Value = 0.06896551724137931
$ws.Cells.Item(4, 15).Value = 0.03448275862068965
$ws.Cells.Item(4, 16).Value = 0.7241379310344828
$ws.Cells.Item(4, 19).Value = 0.1724137931034483
$ws.Cells.Item(5, 19).Value = 1
$ws.Cells.Item(6, 2).Value = 0.0625
$ws.Cells.Item(6, 4).Value = 0.005208333333333333
$ws.Cells.Item(6, 6).Value = 0.05729166666666666
$ws.Cells.Item(6, 10).Value = 0.2552083333333333
$ws.Cells.Item(6, 15).Value = 0.015625
$ws.Cells.Item(6, 18).Value = 0.09895833333333333
$ws.Cells.Item(6, 19).Value = 0.3385416666666667
$ws.Cells.Item(7, 2).Value = 0.1210526315789474
$ws.Cells.Item(7, 4).Value = 0.03684210526315789
$ws.Cells.Item(7, 5).Value = 0.005263157894736842
$ws.Cells.Item(7, 6).Value = 0.03157894736842105
$ws.Cells.Item(7, 10).Value = 0.1157894736842105
$ws.Cells.Item(7, 15).Value = 0.005263157894736842
$ws.Cells.Item(7, 17).Value = 0.1894736842105263
$ws.Cells.Item(7, 18).Value = 0.08947368421052632
$ws.Cells.Item(7, 19).Value = 0.4052631578947368
$ws.Cells.Item(8, 2).Value = 0.08580858085808581
$ws.Cells.Item(8, 4).Value = 0.0231023102310231
$ws.Cells.Item(8, 6).Value = 0.05280528052805281
$ws.Cells.Item(8, 10).Value = 0.1782178217821782
$ws.Cells.Item(8, 15).Value = 0.0231023102310231
$ws.Cells.Item(8, 17).Value = 0.2145214521452145
$ws.Cells.Item(8, 18).Value = 0.1023102310231023
$ws.Cells.Item(8, 19).Value = 0.3201320132013201
$ws.Cells.Item(9, 2).Value = 0.08450704225352113
$ws.Cells.Item(9, 6).Value = 0.04225352112676056
$ws.Cells.Item(9, 10).Value = 0.1126760563380282
$ws.Cells.Item(9, 15).Value = 0.0352112676056338
$ws.Cells.Item(9, 17).Value = 0.1690140845070423
$ws.Cells.Item(9, 18).Value = 0.09154929577464789
$ws.Cells.Item(9, 19).Value = 0.4647887323943662
$ws.Cells.Item(10, 2).Value = 0.09310986964618249
$ws.Cells.Item(10, 4).Value = 0.01675977653631285
$ws.Cells.Item(10, 6).Value = 0.08379888268156424
$ws.Cells.Item(10, 10).Value = 0.1378026070763501
$ws.Cells.Item(10, 15).Value = 0.01210428305400372
$ws.Cells.Item(10, 17).Value = 0.1815642458100559
$ws.Cells.Item(10, 18).Value = 0.09683426443202979
$ws.Cells.Item(10, 19).Value = 0.3780260707635009
$ws.Cells.Item(11, 7).Value = 0.1533333333333333
$ws.Cells.Item(11, 10).Value = 0.08333333333333333
$ws.Cells.Item(11, 11).Value = 0.21
$ws.Cells.Item(11, 12).Value = 0.5266666666666666
$ws.Cells.Item(11, 19).Value = 0.02666666666666667
$ws.Cells.Item(12, 7).Value = 0.7195121951219512
$ws.Cells.Item(12, 10).Value = 0.1890243902439024
$ws.Cells.Item(12, 11).Value = 0.01219512195121951
$ws.Cells.Item(12, 12).Value = 0.03658536585365853
$ws.Cells.Item(12, 19).Value = 0.0426829268292683
$ws.Cells.Item(13, 6).Value = 0.02325581395348837
$ws.Cells.Item(13, 7).Value = 0.7209302325581395
$ws.Cells.Item(13, 10).Value = 0.2093023255813954
$ws.Cells.Item(13, 19).Value = 0.04651162790697674
$ws.Cells.Item(15, 6).Value = 0.03048780487804878
$ws.Cells.Item(15, 8).Value = 0.1463414634146341
$ws.Cells.Item(15, 9).Value = 0.06707317073170732
$ws.Cells.Item(15, 10).Value = 0.3597560975609756
$ws.Cells.Item(15, 11).Value = 0.0975609756097561
$ws.Cells.Item(15, 13).Value = 0.01219512195121951
$ws.Cells.Item(15, 15).Value = 0.06097560975609756
$ws.Cells.Item(15, 19).Value = 0.225609756097561
$ws.Cells.Item(16, 6).Value = 0.006896551724137931
$ws.Cells.Item(16, 8).Value = 0.1724137931034483
$ws.Cells.Item(16, 9).Value = 0.07586206896551724
$ws.Cells.Item(16, 10).Value = 0.3586206896551724
$ws.Cells.Item(16, 11).Value = 0.1655172413793103
$ws.Cells.Item(16, 13).Value = 0.04827586206896552
$ws.Cells.Item(16, 15).Value = 0.06206896551724138
$ws.Cells.Item(16, 19).Value = 0.1103448275862069
$ws.Cells.Item(17, 6).Value = 0.005747126436781609
$ws.Cells.Item(17, 8).Value = 0.1408045977011494
$ws.Cells.Item(17, 9).Value = 0.09770114942528736
$ws.Cells.Item(17, 10).Value = 0.4482758620689655
$ws.Cells.Item(17, 11).Value = 0.1206896551724138
$ws.Cells.Item(17, 13).Value = 0.01149425287356322
$ws.Cells.Item(17, 15).Value = 0.04597701149425287
$ws.Cells.Item(17, 19).Value = 0.1293103448275862
$ws.Cells.Item(18, 6).Value = 0.01639344262295082
$ws.Cells.Item(18, 8).Value = 0.180327868852459
$ws.Cells.Item(18, 9).Value = 0.0546448087431694
$ws.Cells.Item(18, 10).Value = 0.453551912568306
$ws.Cells.Item(18, 11).Value = 0.09289617486338798
$ws.Cells.Item(18, 13).Value = 0.04371584699453552
$ws.Cells.Item(18, 14).Value = 0.00546448087431694
$ws.Cells.Item(18, 15).Value = 0.07103825136612021
$ws.Cells.Item(18, 19).Value = 0.08196721311475409
$ws.Cells.Item(19, 6).Value = 0.01859099804305284
$ws.Cells.Item(19, 8).Value = 0.1692759295499021
$ws.Cells.Item(19, 9).Value = 0.07827788649706457
$ws.Cells.Item(19, 10).Value = 0.3679060665362035
$ws.Cells.Item(19, 11).Value = 0.1281800391389432
$ws.Cells.Item(19, 13).Value = 0.02250489236790607
$ws.Cells.Item(19, 15).Value = 0.0675146771037182
$ws.Cells.Item(19, 19).Value = 0.1477495107632094
